$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename parameter labels: strip unit suffixes / fix Engine Type ---
# (order matters for shared-string table append order - matches target)
$ws.Range("B3").Value  = "Tf"
$ws.Range("B4").Value  = "Tr"
$ws.Range("B5").Value  = "W"
$ws.Range("B6").Value  = "xf"
$ws.Range("B7").Value  = "D_wheel"
$ws.Range("B15").Value = "m_p"
$ws.Range("B16").Value = "m_v"
$ws.Range("B17").Value = "h_p"
$ws.Range("B18").Value = "h_v"
$ws.Range("B22").Value = "S"
$ws.Range("B23").Value = "h_aero"
$ws.Range("B24").Value = "m_aero"
$ws.Range("B39").Value = "h"
$ws.Range("B40").Value = "m_t"
$ws.Range("B25").Value = "Engine_Type"

# --- Data value changes ---
$ws.Range("C19").Value = 2

# --- New values (points per épreuve) ---
$ws.Range("C30").Value = 1111
$ws.Range("D30").Value = 1111

$ws.Range("C33").Value = 1111
$ws.Range("D33").Value = 1111

$ws.Range("C37").Value = 1111
$ws.Range("D37").Value = 1111

$ws.Range("C38").Value = 1111
$ws.Range("D38").Value = 1111

# --- Selection change ---
$null = $ws.Range("G17").Select()
